$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Rename header columns in row 1: "_old" -> "_FV2310", "_new" -> "_FV2404" ---
$baseNames = @("Segmentname", "Segmentgruppe", "Segment", "Datenelement", "Segment ID", "Code", "Qualifier", "Beschreibung", "Bedingungsausdruck", "Bedingung")

for ($i = 0; $i -lt $baseNames.Count; $i++) {
    # Columns A-J (1-10): "<name>_old" -> "<name>_FV2310"
    $ws.Cells.Item(1, $i + 1).Value = ($baseNames[$i] + "_FV2310")
    # Columns L-U (12-21): "<name>_new" -> "<name>_FV2404"
    $ws.Cells.Item(1, $i + 12).Value = ($baseNames[$i] + "_FV2404")
}

# Column K (11) stays "diff" - unchanged.

# --- 2. Freeze the header row (split/freeze at row 2, pane anchored top-left A2) ---
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true

# --- 3. Turn the used range A1:U67 into an Excel Table (ListObject) named Table1 ---
$tbl = $ws.ListObjects.Add(1, $ws.Range("A1:U67"), $null, 1)
$tbl.Name = "Table1"
$tbl.TableStyle = ""
